$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header row (row 1): replace the numeric sequence 1..40 in B1:AO1 with
#    text labels "var1".."var40" (pulled from the re-queried database).
# ---------------------------------------------------------------------------

# F1 currently carries the red/white "highlight" style (border+fill+white
# font). In the rebuilt header it should look like every other header cell,
# so copy the plain look from B1 onto it before writing the new values.
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

for ($i = 1; $i -le 40; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "var$i"
}

# ---------------------------------------------------------------------------
# 2) Small KPI block below the data grid (rows 6-10) gets rebuilt: a thin
#    spacer row is added, the old formulas are wiped, and the surviving
#    placeholder cells are relocated with their original number formats.
# ---------------------------------------------------------------------------

# Stash T7's current format (#,##0 style) in an unused cell so it can be
# re-applied after T7 and T8 swap roles below.
$holdT7 = $ws.Range("ZZ1")
$ws.Range("T7").Copy()
$holdT7.PasteSpecial(-4122)

# New T10 inherits the old T7 formatting (#,##0 style).
$ws.Range("T7").Copy()
$ws.Range("T10").PasteSpecial(-4122)

# T7 becomes the old T8 formatting (percentage style).
$ws.Range("T8").Copy()
$ws.Range("T7").PasteSpecial(-4122)

# T8 becomes the old T7 formatting (#,##0 style), restored from the stash.
$holdT7.Copy()
$ws.Range("T8").PasteSpecial(-4122)

# New B10 / F10 inherit the old B7 / F7 formatting.
$ws.Range("B7").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("F7").Copy()
$ws.Range("F10").PasteSpecial(-4122)

# New C7 inherits the old C8 formatting.
$ws.Range("C8").Copy()
$ws.Range("C7").PasteSpecial(-4122)

$holdT7.Clear()

# Wipe the stale formulas/values from the cells that are kept in place.
$ws.Range("T7").ClearContents()
$ws.Range("T8").ClearContents()

# Drop the cells that no longer have any role at their old address.
$ws.Range("B7").Clear()
$ws.Range("F7").Clear()
$ws.Range("C8").Clear()
$ws.Range("T9").Clear()

# Thin spacer row separating the data grid from the KPI placeholders.
$ws.Rows.Item(6).RowHeight = 1.5

Write-Output "edit applied"
